$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDA")
Write-Host $ws.Name
